# Sync attendance_reports, modules_schedules, and assets from main repo - 2026-01-17 12:29:19
#
# This script updates the "Recorded By" column (G) on every recorded
# attendance row from "System, dnasr281@gmail.com" to
# "dnasr281@gmail.com, System", and refreshes a handful of recomputed
# attendance counts / percentages elsewhere on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the "Recorded By" ordering on every row that has it -------------
$lastRow = $ws.UsedRange.Rows.Count
$oldVal = "System, dnasr281@gmail.com"
$newVal = "dnasr281@gmail.com, System"
$swapped = 0

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value()
    if ($val -eq $oldVal) {
        $cell.Value = $newVal
        $swapped = $swapped + 1
    }
}
Write-Host "Recorded By swapped on $swapped rows"

# --- Updated attendance counts (students present / total) -----------------
# These are plain "x/y" strings (not date-like to Excel), so a direct
# .Value assignment keeps them as text with their original formatting.
$ws.Range("H178").Value = "28/30"
$ws.Range("H195").Value = "27/27"
$ws.Range("H216").Value = "23/29"
$ws.Range("H237").Value = "24/29"

# --- Recomputed attendance percentages -------------------------------------
# These cells store the percentage as literal text (e.g. "81.1%"), not a
# numeric percent value. Force the Text number format first so Excel's
# auto-detection doesn't silently convert the string into a numeric
# percentage (which would change the cell's stored type/format).
$pctCells = @("L10", "S23", "S24", "S25", "S26")
foreach ($ref in $pctCells) {
    $r = $ws.Range($ref)
    $r.NumberFormat = "@"
}

$ws.Range("L10").Value = "81.1%"
$ws.Range("S23").Value = "85.0%"
$ws.Range("S24").Value = "73.5%"
$ws.Range("S25").Value = "76.1%"
$ws.Range("S26").Value = "75.3%"

Write-Host "Done"
